# One more BS32 run to go
# Fill in the missing "Av. integr. Calls" (column D) results for the BS32
# (rows 10/11/13), BS54 (row 26) and DP54 (row 39) tolerance tables, let the
# existing "=n*D.." formulas in column E recompute, and refresh the
# already-collected BS32 call-count average in G33 with the latest run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D10").Value = 1991
$ws.Range("D11").Value = 4293
$ws.Range("D13").Value = 19929
$ws.Range("D26").Value = 613
$ws.Range("D39").Value = 799

$ws.Range("G33").Formula = "=(4223+4269+4457+4221)/4"

# The added numbers are wider than the previous content, so the sheet's
# (non-custom) "optimal" column widths shrink slightly once recalculated.
# Reproduce the resulting widths as closely as this engine's column-width
# granularity allows.
$ws.Columns("A:C").ColumnWidth = 7.5
$ws.Columns("D").ColumnWidth = 15.83
$ws.Columns("E").ColumnWidth = 15.17
$ws.Columns("F:AMK").ColumnWidth = 7.5

$ws.Range("D12").Select()
